# Apply crypto price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '65.731.85'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').Value = '3.419.16'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.77'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.93'
$ws.Range('E6').Value = '  -3.70%  '
$ws.Range('D7').Value = '3.420.44'
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('E10').Value = '  -4.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.80'
$ws.Range('E11').Value = '  +4.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.405'
$ws.Range('E12').Value = '  -4.15%  '
$ws.Range('D13').Value = '4.007.16'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000202'
$ws.Range('E14').Value = '  -4.88%  '
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('D16').Value = '3.431.34'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '65.846.83'
$ws.Range('E17').Value = '  -1.43%  '
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.42'
$ws.Range('E19').Value = '  +4.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.11'
$ws.Range('E20').Value = '  -4.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.84'
$ws.Range('E21').Value = '  -2.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '416.78'
$ws.Range('E22').Value = '  -4.95%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.582'
$ws.Range('E23').Value = '  -4.92%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.66'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '3.551.74'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000112'
$ws.Range('E27').Value = '  -6.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.26'
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.79'
$ws.Range('E29').Value = '  -6.53%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.44'
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.161'
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.46'
$ws.Range('E33').Value = '  -8.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.62'
$ws.Range('E34').Value = '  -2.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.71'
$ws.Range('E36').Value = '  -4.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.54'
$ws.Range('E37').Value = '  -9.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.63'
$ws.Range('E38').Value = '  -3.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '173.54'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0865'
$ws.Range('E41').Value = '  -2.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.09'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.870'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.92'
$ws.Range('E44').Value = '  -11.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '45.51'
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.79'
$ws.Range('E46').Value = '  -7.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.17'
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.08'
$ws.Range('E48').Value = '  -5.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.28'
$ws.Range('E49').Value = '  -6.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.924'
$ws.Range('E50').Value = '  -6.24%  '
$ws.Range('E51').Value = '  -4.66%  '
